# Update countries & provincias Spain
# Applies the data refresh captured in the commit diff:
#  - Updates the "Datos actualizados..." timestamp string (cell A1)
#  - Updates numeric stats for several country rows
#  - Because Belgica/Marruecos, Birmania/Uganda and Islas Malvinas/Montserrat
#    now have swapped relative totals, their row order (column A labels) swaps too

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header / timestamp -----------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 22 de Septiembre de 2020 a las 05:20"

# --- Helper table of row updates ---------------------------------------
# Each entry: row number, B..H values (country name in A left untouched
# unless noted further below)
$rows = @(
    @{ Row = 21;  B = 306886; C = 582;  D = 293159; E = 7303;  F = 0; G = 4;  H = 6424 },
    @{ Row = 30;  B = 130986; C = 310;  D = 90240;  E = 33092; F = 0; G = 37; H = 7654 },
    @{ Row = 35;  B = 107374; C = 67;   D = 102064; E = 3639;  F = 0; G = 0;  H = 1671 },
    @{ Row = 37;  B = 103392; C = 1097; D = 18977;  E = 74465; F = 0; G = 2;  H = 9950 },
    @{ Row = 38;  B = 103119; C = 0;    D = 84158;  E = 17106; F = 0; G = 0;  H = 1855 },
    @{ Row = 50;  B = 72075;  C = 459;  D = 22611;  E = 47260; F = 0; G = 20; H = 2204 },
    @{ Row = 78;  B = 26942;  C = 30;   D = 24157;  E = 1931;  F = 0; G = 3;  H = 854 },
    @{ Row = 112; B = 6471;   C = 320;  D = 1445;   E = 4926;  F = 0; G = 2;  H = 100 },
    @{ Row = 113; B = 6468;   C = 0;    D = 2731;   E = 3674;  F = 0; G = 0;  H = 63 },
    @{ Row = 159; B = 1635;   C = 8;    D = 967;    E = 647;   F = 0; G = 0;  H = 21 },
    @{ Row = 173; B = 591;    C = 7;    D = 505;    E = 66;    F = 0; G = 0;  H = 20 },
    @{ Row = 185; B = 313;    C = 1;    D = 302;    E = 11;    F = 0; G = 0;  H = 0 },
    @{ Row = 214; B = 13;     C = 0;    D = 13;     E = 0;     F = 0; G = 0;  H = 0 },
    @{ Row = 215; B = 13;     C = 0;    D = 12;     E = 0;     F = 0; G = 0;  H = 1 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
}

# --- Row-order swaps (country names in column A) ------------------------
# Updated totals flip the relative ranking of these country pairs, so the
# labels in column A swap while the (already updated) numeric data stays
# attached to its row position.
$ws.Range("A37").Value = "Belgica"
$ws.Range("A38").Value = "Marruecos"

$ws.Range("A112").Value = "Birmania"
$ws.Range("A113").Value = "Uganda"

$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("A215").Value = "Montserrat"
